$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add "Diff x" / "Diff y" headers in E1 / F1 ---------------------------
# Write the plain text first, then paste-in just the cell formatting from an
# existing header cell (D1 = "Predicted yi") so the new headers inherit the
# same base font / fill / alignment as the rest of row 1.
$ws.Range("E1").Value = "Diff x"
$ws.Range("F1").Value = "Diff y"

$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Italicize / re-font the last character ("x" / "y") to match the other
# headers' styling of the trailing run.
$ex = $ws.Range("E1").Characters(6,1)
$ex.Font.Italic = $true
$ex.Font.Name = "Aileron Heavy"

$fy = $ws.Range("F1").Characters(6,1)
$fy.Font.Italic = $true
$fy.Font.Name = "Aileron Heavy"

# --- 2. Diff formulas for rows 2-21 ------------------------------------------
$ws.Range("E2").Formula = "=ABS(A2-C2)"
$ws.Range("F2").Formula = "=ABS(B2-D2)"

$ws.Range("E3:E21").Formula = "=ABS(A3-C3)"
$ws.Range("F3:F21").Formula = "=ABS(B3-D3)"

# --- 3. Average row (22) ------------------------------------------------------
$ws.Range("E22").Formula = "=AVERAGE(E2:E21)"
$ws.Range("F22").Formula = "=AVERAGE(F2:F21)"

$avgRange = $ws.Range("E22:F22")
$avgRange.Font.Bold = $true
$avgRange.Font.Color = 255

# --- 4. Remove the orange highlight fill from A2 / C2 -------------------------
# A2 / C2 keep their current fonts but should no longer have the orange fill.
# Re-use the formatting of sibling cells that already have the same font with
# no fill (B2 for A2's font; C3 for C2's font) instead of poking Interior
# directly, which would otherwise fork brand new font/fill table entries.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 5. Selection + autofit -----------------------------------------------
$ws.Range("J17").Select()
